# Actualiza la planilla de métricas: completa la fila 1 y 2 de la tabla
# "Desarrollo y correctivos" (Incrementos 1 y 2) con los datos reales de
# la tarea "Especie".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Incremento 1: "Especie: atrubutos, getters y const"
$ws.Range("C18").Value2 = "Especie: atrubutos, getters y const"
$ws.Range("F18").Value2 = 20
$ws.Range("G18").Value2 = 0.003472222222222222
$ws.Range("H18").Value2 = 0.6909722222222222
$ws.Range("I18").Value2 = 0.6923611111111111
$ws.Range("K18").Value2 = 0
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = 17

# Incremento 2: "Especie: implementa Comparable"
$ws.Range("C19").Value2 = "Especie: implementa Comparable"
$ws.Range("F19").Value2 = 3
$ws.Range("G19").Value2 = 0.0020833333333333333
$ws.Range("H19").Value2 = 0.6923611111111111
$ws.Range("I19").Value2 = 0.6930555555555555
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 0
$ws.Range("M19").Value2 = 3

# Deja seleccionado el mismo rango que quedó activo en el archivo original
$ws.Range("C20:E20").Select()
